$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2-51 Price (D) and Volume(1h) (E) updates ---
$ws.Range("D2").Value = "67.603.27"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "3.321.13"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'578.69"
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("D6").Value = "'174.72"
$ws.Range("E6").Value = "  -4.41%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.68%  "
$ws.Range("D9").Value = "3.316.73"
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("D12").Value = "'45.29"
$ws.Range("E12").Value = "  -2.39%  "
$ws.Range("E13").Value = "  -1.85%  "
$ws.Range("D14").Value = "'656.65"
$ws.Range("E14").Value = "  +3.34%  "
$ws.Range("D15").Value = "3.864.81"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("E16").Value = "  -1.08%  "
$ws.Range("D17").Value = "67.619.13"
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("D19").Value = "3.328.46"
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("E20").Value = "  -2.17%  "
$ws.Range("D21").Value = "'10.94"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("E22").Value = "  -1.99%  "
$ws.Range("D23").Value = "'5.34"
$ws.Range("E23").Value = "  +6.00%  "
$ws.Range("D24").Value = "'17.00"
$ws.Range("E24").Value = "  -3.86%  "
$ws.Range("D25").Value = "'98.50"
$ws.Range("E25").Value = "  +1.50%  "
$ws.Range("E26").Value = "  -4.16%  "
$ws.Range("E27").Value = "  -4.21%  "
$ws.Range("E28").Value = "  -3.89%  "
$ws.Range("D29").Value = "'33.39"
$ws.Range("E29").Value = "  +2.50%  "
$ws.Range("D30").Value = "'8.40"
$ws.Range("E30").Value = "  -2.64%  "
$ws.Range("D31").Value = "'7.21"
$ws.Range("E31").Value = "  +7.12%  "
$ws.Range("D32").Value = "'567.77"
$ws.Range("E32").Value = "  -4.49%  "
$ws.Range("D33").Value = "'10.91"
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("D36").Value = "'56.53"
$ws.Range("E36").Value = "  +1.54%  "
$ws.Range("D37").Value = "3.667.37"
$ws.Range("E37").Value = "  -7.25%  "
$ws.Range("E38").Value = "  -7.39%  "
$ws.Range("D39").Value = "'34.21"
$ws.Range("E39").Value = "  +4.87%  "
$ws.Range("E40").Value = "  +0.88%  "
$ws.Range("E41").Value = "  -2.96%  "
$ws.Range("E44").Value = "  -1.79%  "
$ws.Range("D45").Value = "0.0₃0660"
$ws.Range("E45").Value = "  -3.72%  "
$ws.Range("E46").Value = "  -2.12%  "
$ws.Range("E47").Value = "  +1.68%  "
$ws.Range("E48").Value = "  -1.00%  "
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("E50").Value = "  -1.61%  "
$ws.Range("D51").Value = "'129.28"
$ws.Range("E51").Value = "  -0.91%  "

# --- Rows 42/43: Stacks / ApeXProtocol swap positions ---
$ws.Range("B42").Value = "ApeXProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D42").Value = "'3.36"
$ws.Range("E42").Value = "  -1.07%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'3.10"
$ws.Range("E43").Value = "  -4.62%  "
